$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates ---
$ws.Range("B3").Value = "No"
$ws.Range("L3").Value = 28

# M3: slot time (12:30 PM) with h:mm time number format
$ws.Range("M3").Value = 0.52083333333333337
$ws.Range("M3").NumberFormat = "h:mm"

# N3: Session = PM (copy format from L3 which already has the centered bordered style)
$ws.Range("L3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null
$ws.Range("N3").Value = "PM"

# --- Row 4 updates ---
$ws.Range("D4").Value = "test@gmail.com"
$ws.Range("I4").Value = "Yes"

# J4 / K4 were empty with a plain bordered style; add values + center alignment
$ws.Range("J4").Value = "Dr. Emily Larson"
$ws.Range("J4").HorizontalAlignment = -4108
$ws.Range("K4").Value = 649
$ws.Range("K4").HorizontalAlignment = -4108

$ws.Range("L4").Value = 30

# M4: slot time (11:00 AM) with h:mm time number format
$ws.Range("M4").Value = 0.45833333333333331
$ws.Range("M4").NumberFormat = "h:mm"

# N4: Session = AM
$ws.Range("L4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").Value = "AM"

# --- Header row 2: new "Session" column, copy header style from the previous header cell ---
$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null
$ws.Range("N2").Value = "Session"

# --- Selection update ---
$ws.Range("H12").Select() | Out-Null
